# Insert a new record row into the daily "Puerro" price log at row 52,
# pushing the existing rows 52..138 down to 53..139 (weekly refresh of the
# "Fruta / hortaliza, semanal" subset).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 52..138 down one row, creating a blank row 52.
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row 52 with the new daily record.
$ws.Range("A52").Value = 9
$ws.Range("B52").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C52").Value = "Metropolitana"
$ws.Range("D52").Value = 45210
$ws.Range("E52").Value = 13
$ws.Range("F52").Value = 100112005
$ws.Range("G52").Value = "Puerro"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 70
$ws.Range("K52").Value = 8000
$ws.Range("L52").Value = 8000
$ws.Range("M52").Value = 8000
$ws.Range("N52").Value = "$/paquete 20 unidades"
$ws.Range("O52").Value = "Provincia de Chacabuco"
$ws.Range("P52").Value = 400
$ws.Range("Q52").Value = 20
$ws.Range("R52").Value = "Hortaliza"
